# Generate Report for Handback
# Refresh the "Latest HO Xliff Generate Date" / "Correspond Handoff Datetime" /
# "Correspond Handback DateTime" timestamps for the 2c29110a... file (and the
# matching handoff time on the de-de sheet) to reflect a newly generated
# handback report.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: row 3 corresponds to 2c29110a-a729-4e33-bcb9-89a01b98271b.md
$wsOverview.Range("G3").Value = "2016-08-23 14:58:36"

# zh-cn sheet: row 3 corresponds to 2c29110a-a729-4e33-bcb9-89a01b98271b.md
$wsZhCn.Range("H3").Value = "2016-08-23 14:58:31"
$wsZhCn.Range("K3").Value = "2016-08-23 14:58:49"

# de-de sheet: row 3 corresponds to 2c29110a-a729-4e33-bcb9-89a01b98271b.md
$wsDeDe.Range("H3").Value = "2016-08-23 14:58:36"
$wsDeDe.Range("K3").Value = "2016-08-23 14:58:57"
